# Realestate Update resale numbers 2025-02-28 22:44
# Appends the new daily resale-number snapshot row (row 94) to the
# "CityResaleNum" sheet, mirroring the existing data rows: columns A-D are
# text (date/time/weekday/week-number-as-text), columns E-T are numeric
# city values (-1 = not reported).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 94
$textRange = "A" + $row + ":D" + $row

# --- text columns -----------------------------------------------------
# A leading apostrophe forces these to be stored as genuine text (matching
# how every prior row stores Date/Time/Weekday/Week as text, not as an
# auto-converted date serial / number), then the style is put back to
# Normal so no visible formatting change is introduced.
$ws.Range("A$row").Value = "'2025-02-28"
$ws.Range("B$row").Value = "'22:44:30"
$ws.Range("C$row").Value = "'Friday"
$ws.Range("D$row").Value = "'08"
$ws.Range($textRange).Style = "Normal"

# --- numeric city columns ----------------------------------------------
$ws.Range("E$row").Value = 131933
$ws.Range("F$row").Value = 142381
$ws.Range("G$row").Value = 173357
$ws.Range("H$row").Value = 160691
$ws.Range("I$row").Value = -1
$ws.Range("J$row").Value = 147433
$ws.Range("K$row").Value = -1
$ws.Range("L$row").Value = -1
$ws.Range("M$row").Value = 194828
$ws.Range("N$row").Value = 115831
$ws.Range("O$row").Value = 47093
$ws.Range("P$row").Value = 29755
$ws.Range("Q$row").Value = 70577
$ws.Range("R$row").Value = -1
$ws.Range("S$row").Value = 51652
$ws.Range("T$row").Value = -1
